# Add I0 (column I) and IF (column J) headers + data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers, copying the style used by the existing header row (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row data: row number, I value, J value
$data = @(
    @(2, 8, 8),
    @(3, 7, 7),
    @(4, 6, 6),
    @(5, 7, 8),
    @(6, 7, 7),
    @(7, 9, 9),
    @(8, 8, 8),
    @(9, 5, 6),
    @(10, 8, 8),
    @(11, 9, 9),
    @(12, 7, 7),
    @(13, 5, 6),
    @(14, 6, 6),
    @(15, 9, 9),
    @(16, 7, 7),
    @(17, 6, 7),
    @(18, 9, 9),
    @(19, 7, 7),
    @(20, 8, 8),
    @(21, 9, 9),
    @(22, 9, 9),
    @(23, 8, 8),
    @(24, 9, 9),
    @(25, 7, 7),
    @(26, 9, 9),
    @(27, 7, 7),
    @(28, 9, 9),
    @(29, 9, 9),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 8, 9),
    @(34, 9, 9),
    @(35, 8, 8),
    @(36, 9, 9),
    @(37, 9, 9),
    @(38, 9, 9),
    @(39, 8, 8),
    @(40, 9, 9),
    @(41, 8, 8),
    @(42, 9, 9),
    @(43, 8, 8),
    @(44, 9, 9),
    @(45, 9, 9),
    @(46, 8, 8),
    @(47, 8, 8),
    @(48, 8, 8),
    @(49, 8, 8),
    @(50, 9, 9),
    @(51, 8, 8),
    @(52, 7, 7),
    @(53, 7, 7),
    @(54, 6, 6),
    @(55, 5, 5),
    @(56, 6, 7),
    @(57, 5, 5),
    @(58, 4, 4),
    @(59, 3, 3),
    @(60, 5, 5)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
